$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows whose Target cluster is "Resolving-Mac" (bottom-up so row indices stay valid)
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(6).Delete()

# Refresh the recomputed TPM-derived numeric columns (G,H,I,J,M,N,O,P,Q,R,S,T) for the remaining rows
# row 2
$ws.Range("G2").Value = [double]"0.026992"
$ws.Range("H2").Value = [double]"0.08097599999999999"
$ws.Range("I2").Value = [double]"0.004182906599909731"
$ws.Range("J2").Value = [double]"0.00420788870005516"
$ws.Range("M2").Value = [double]"31.61349966666667"
$ws.Range("N2").Value = [double]"94.84049899999999"
$ws.Range("O2").Value = [double]"0.9864449253320818"
$ws.Range("P2").Value = [double]"0.9897213119532253"
$ws.Range("Q2").Value = [double]"0.8533115830026666"
$ws.Range("R2").Value = [double]"7.679804247023998"
$ws.Range("S2").Value = [double]"0.004126206988619026"
$ws.Range("T2").Value = [double]"0.004164637124771746"

# row 3
$ws.Range("G3").Value = [double]"0.026992"
$ws.Range("H3").Value = [double]"0.08097599999999999"
$ws.Range("I3").Value = [double]"0.004182906599909731"
$ws.Range("J3").Value = [double]"0.00420788870005516"
$ws.Range("M3").Value = [double]"0.064626"
$ws.Range("N3").Value = [double]"0.193878"
$ws.Range("O3").Value = [double]"0.002016543262109296"
$ws.Range("P3").Value = [double]"0.002023241026166126"
$ws.Range("Q3").Value = [double]"0.001744384992"
$ws.Range("R3").Value = [double]"0.015699464928"
$ws.Range("S3").Value = [double]"8.43501212008047E-06"
$ws.Range("T3").Value = [double]"8.513573051492446E-06"

# row 4
$ws.Range("G4").Value = [double]"0.026992"
$ws.Range("H4").Value = [double]"0.08097599999999999"
$ws.Range("I4").Value = [double]"0.004182906599909731"
$ws.Range("J4").Value = [double]"0.00420788870005516"
$ws.Range("M4").Value = [double]"0.05151033333333333"
$ws.Range("N4").Value = [double]"0.154531"
$ws.Range("O4").Value = [double]"0.001607291424694971"
$ws.Range("P4").Value = [double]"0.001612629896194914"
$ws.Range("Q4").Value = [double]"0.001390366917333333"
$ws.Range("R4").Value = [double]"0.012513302256"
$ws.Range("S4").Value = [double]"6.723149908334906E-06"
$ws.Range("T4").Value = [double]"6.785767117569706E-06"

# row 5
$ws.Range("G5").Value = [double]"0.026992"
$ws.Range("H5").Value = [double]"0.08097599999999999"
$ws.Range("I5").Value = [double]"0.004182906599909731"
$ws.Range("J5").Value = [double]"0.00420788870005516"
$ws.Range("M5").Value = [double]"0.3182755"
$ws.Range("N5").Value = [double]"0.636551"
$ws.Range("O5").Value = [double]"0.009931239981113902"
$ws.Range("P5").Value = [double]"0.00664281712441367"
$ws.Range("Q5").Value = [double]"0.008590892296"
$ws.Range("R5").Value = [double]"0.051545353776"
$ws.Range("S5").Value = [double]"4.154144926228873E-05"
$ws.Range("T5").Value = [double]"2.79522351143532E-05"

# row 6
$ws.Range("G6").Value = [double]"6.300519666666666"
$ws.Range("H6").Value = [double]"18.901559"
$ws.Range("I6").Value = [double]"0.976381346197431"
$ws.Range("J6").Value = [double]"0.9822127115383066"
$ws.Range("M6").Value = [double]"31.61349966666667"
$ws.Range("N6").Value = [double]"94.84049899999999"
$ws.Range("O6").Value = [double]"0.9864449253320818"
$ws.Range("P6").Value = [double]"0.9897213119532253"
$ws.Range("Q6").Value = [double]"199.1814763819934"
$ws.Range("R6").Value = [double]"1792.633287437941"
$ws.Range("S6").Value = [double]"0.9631464241453623"
$ws.Range("T6").Value = [double]"0.9721168534808277"

# row 7
$ws.Range("G7").Value = [double]"6.300519666666666"
$ws.Range("H7").Value = [double]"18.901559"
$ws.Range("I7").Value = [double]"0.976381346197431"
$ws.Range("J7").Value = [double]"0.9822127115383066"
$ws.Range("M7").Value = [double]"0.064626"
$ws.Range("N7").Value = [double]"0.193878"
$ws.Range("O7").Value = [double]"0.002016543262109296"
$ws.Range("P7").Value = [double]"0.002023241026166126"
$ws.Range("Q7").Value = [double]"0.407177383978"
$ws.Range("R7").Value = [double]"3.664596455802"
$ws.Range("S7").Value = [double]"0.001968915224923633"
$ws.Range("T7").Value = [double]"0.001987253054406176"

# row 8
$ws.Range("G8").Value = [double]"6.300519666666666"
$ws.Range("H8").Value = [double]"18.901559"
$ws.Range("I8").Value = [double]"0.976381346197431"
$ws.Range("J8").Value = [double]"0.9822127115383066"
$ws.Range("M8").Value = [double]"0.05151033333333333"
$ws.Range("N8").Value = [double]"0.154531"
$ws.Range("O8").Value = [double]"0.001607291424694971"
$ws.Range("P8").Value = [double]"0.001612629896194914"
$ws.Range("Q8").Value = [double]"0.3245418682032222"
$ws.Range("R8").Value = [double]"2.920876813829"
$ws.Range("S8").Value = [double]"0.001569329364975262"
$ws.Range("T8").Value = [double]"0.001583945583049345"

# row 9
$ws.Range("G9").Value = [double]"6.300519666666666"
$ws.Range("H9").Value = [double]"18.901559"
$ws.Range("I9").Value = [double]"0.976381346197431"
$ws.Range("J9").Value = [double]"0.9822127115383066"
$ws.Range("M9").Value = [double]"0.3182755"
$ws.Range("N9").Value = [double]"0.636551"
$ws.Range("O9").Value = [double]"0.009931239981113902"
$ws.Range("P9").Value = [double]"0.00664281712441367"
$ws.Range("Q9").Value = [double]"2.005301047168166"
$ws.Range("R9").Value = [double]"12.031806283009"
$ws.Range("S9").Value = [double]"0.009696677462169741"
$ws.Range("T9").Value = [double]"0.006524659420023447"

# row 10
$ws.Range("G10").Value = [double]"0.1149325"
$ws.Range("H10").Value = [double]"0.229865"
$ws.Range("I10").Value = [double]"0.01781090370458377"
$ws.Range("J10").Value = [double]"0.01194485200600399"
$ws.Range("M10").Value = [double]"31.61349966666667"
$ws.Range("N10").Value = [double]"94.84049899999999"
$ws.Range("O10").Value = [double]"0.9864449253320818"
$ws.Range("P10").Value = [double]"0.9897213119532253"
$ws.Range("Q10").Value = [double]"3.633418550439167"
$ws.Range("R10").Value = [double]"21.800511302635"
$ws.Range("S10").Value = [double]"0.01756947557496503"
$ws.Range("T10").Value = [double]"0.01182207459846939"

# row 11
$ws.Range("G11").Value = [double]"0.1149325"
$ws.Range("H11").Value = [double]"0.229865"
$ws.Range("I11").Value = [double]"0.01781090370458377"
$ws.Range("J11").Value = [double]"0.01194485200600399"
$ws.Range("M11").Value = [double]"0.064626"
$ws.Range("N11").Value = [double]"0.193878"
$ws.Range("O11").Value = [double]"0.002016543262109296"
$ws.Range("P11").Value = [double]"0.002023241026166126"
$ws.Range("Q11").Value = [double]"0.007427627745"
$ws.Range("R11").Value = [double]"0.04456576647"
$ws.Range("S11").Value = [double]"3.591645785755589E-05"
$ws.Range("T11").Value = [double]"2.416731463003002E-05"

# row 12
$ws.Range("G12").Value = [double]"0.1149325"
$ws.Range("H12").Value = [double]"0.229865"
$ws.Range("I12").Value = [double]"0.01781090370458377"
$ws.Range("J12").Value = [double]"0.01194485200600399"
$ws.Range("M12").Value = [double]"0.05151033333333333"
$ws.Range("N12").Value = [double]"0.154531"
$ws.Range("O12").Value = [double]"0.001607291424694971"
$ws.Range("P12").Value = [double]"0.001612629896194914"
$ws.Range("Q12").Value = [double]"0.005920211385833333"
$ws.Range("R12").Value = [double]"0.035521268315"
$ws.Range("S12").Value = [double]"2.862731279044538E-05"
$ws.Range("T12").Value = [double]"1.926262545050583E-05"

# row 13
$ws.Range("G13").Value = [double]"0.1149325"
$ws.Range("H13").Value = [double]"0.229865"
$ws.Range("I13").Value = [double]"0.01781090370458377"
$ws.Range("J13").Value = [double]"0.01194485200600399"
$ws.Range("M13").Value = [double]"0.3182755"
$ws.Range("N13").Value = [double]"0.636551"
$ws.Range("O13").Value = [double]"0.009931239981113902"
$ws.Range("P13").Value = [double]"0.00664281712441367"
$ws.Range("Q13").Value = [double]"0.03658019890375"
$ws.Range("R13").Value = [double]"0.146320795615"
$ws.Range("S13").Value = [double]"0.000176884358970732"
$ws.Range("T13").Value = [double]"7.934746745407031E-05"

# row 14
$ws.Range("G14").Value = [double]"0.010485"
$ws.Range("H14").Value = [double]"0.031455"
$ws.Range("I14").Value = [double]"0.001624843498075486"
$ws.Range("J14").Value = [double]"0.001634547755634201"
$ws.Range("M14").Value = [double]"31.61349966666667"
$ws.Range("N14").Value = [double]"94.84049899999999"
$ws.Range("O14").Value = [double]"0.9864449253320818"
$ws.Range("P14").Value = [double]"0.9897213119532253"
$ws.Range("Q14").Value = [double]"0.331467544005"
$ws.Range("R14").Value = [double]"2.983207896045"
$ws.Range("S14").Value = [double]"0.001602818623135392"
$ws.Range("T14").Value = [double]"0.001617746749156481"

# row 15
$ws.Range("G15").Value = [double]"0.010485"
$ws.Range("H15").Value = [double]"0.031455"
$ws.Range("I15").Value = [double]"0.001624843498075486"
$ws.Range("J15").Value = [double]"0.001634547755634201"
$ws.Range("M15").Value = [double]"0.064626"
$ws.Range("N15").Value = [double]"0.193878"
$ws.Range("O15").Value = [double]"0.002016543262109296"
$ws.Range("P15").Value = [double]"0.002023241026166126"
$ws.Range("Q15").Value = [double]"0.00067760361"
$ws.Range("R15").Value = [double]"0.006098432489999999"
$ws.Range("S15").Value = [double]"3.27656720802622E-06"
$ws.Range("T15").Value = [double]"3.307084078426878E-06"

# row 16
$ws.Range("G16").Value = [double]"0.010485"
$ws.Range("H16").Value = [double]"0.031455"
$ws.Range("I16").Value = [double]"0.001624843498075486"
$ws.Range("J16").Value = [double]"0.001634547755634201"
$ws.Range("M16").Value = [double]"0.05151033333333333"
$ws.Range("N16").Value = [double]"0.154531"
$ws.Range("O16").Value = [double]"0.001607291424694971"
$ws.Range("P16").Value = [double]"0.001612629896194914"
$ws.Range("Q16").Value = [double]"0.0005400858449999999"
$ws.Range("R16").Value = [double]"0.004860772605"
$ws.Range("S16").Value = [double]"2.611597020928108E-06"
$ws.Range("T16").Value = [double]"2.635920577494012E-06"

# row 17
$ws.Range("G17").Value = [double]"0.010485"
$ws.Range("H17").Value = [double]"0.031455"
$ws.Range("I17").Value = [double]"0.001624843498075486"
$ws.Range("J17").Value = [double]"0.001634547755634201"
$ws.Range("M17").Value = [double]"0.3182755"
$ws.Range("N17").Value = [double]"0.636551"
$ws.Range("O17").Value = [double]"0.009931239981113902"
$ws.Range("P17").Value = [double]"0.00664281712441367"
$ws.Range("Q17").Value = [double]"0.0033371186175"
$ws.Range("R17").Value = [double]"0.020022711705"
$ws.Range("S17").Value = [double]"1.613671071114024E-05"
$ws.Range("T17").Value = [double]"1.08580018217988E-05"

